$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Resolve shapes by their (stable) shape Id rather than positional index.
# ---------------------------------------------------------------------------
function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$shp6  = Get-ShapeById $s 6    # "Initialize" (left column)
$shp7  = Get-ShapeById $s 7    # "Before Change" (left column)
$shp8  = Get-ShapeById $s 8    # "Actions / Formulas" (left column)
$shp9  = Get-ShapeById $s 9    # "After Change" (left column)
$shp10 = Get-ShapeById $s 10   # "Before Save" (left column)
$shp11 = Get-ShapeById $s 11   # "After Save" (left column) -> re-purposed text
$shp12 = Get-ShapeById $s 12   # right-hand down arrow
$shp13 = Get-ShapeById $s 13   # "Initialize" (right column)
$shp14 = Get-ShapeById $s 14   # "On Change" (right column)
$shp17 = Get-ShapeById $s 17   # "Before Save" (right column)

# ---------------------------------------------------------------------------
# A) Re-position existing shapes (left column nudges + right column reflow).
# ---------------------------------------------------------------------------
$shp7.Top  = 120.331730     # 1529862 -> 1528213 EMU
$shp8.Top  = 170.740395     # 2171701 -> 2168403 EMU
$shp9.Top  = 221.149056     # 2800349 -> 2808593 EMU

$shp12.Top    = 56.769293   # 888024  -> 720970  EMU
$shp12.Height = 292.864731  # 2387110 -> 3719382 EMU

$shp13.Left = 560.220185    # 7114797 -> 7114796 EMU
$shp13.Top  = 69.923070     # 1055077 -> 888023  EMU

$shp14.Left = 560.220063    # 7114797 -> 7114795 EMU
$shp14.Top  = 120.331730    # 1696916 -> 1528213 EMU

$shp17.Left = 560.220002    # 7114797 -> 7114794 EMU
$shp17.Top  = 271.557724    # 2409093 -> 3448783 EMU

# ---------------------------------------------------------------------------
# B) Add the four new boxes *before* relabeling $shp11, cloning an existing
#    box so the theme-based line/fill/effect/font style (<p:style>) and the
#    still-unmodified "After Save" run formatting are reproduced exactly.
# ---------------------------------------------------------------------------

# New "After Save" box (accent2 style, like the old shp11) under the left column.
$new15 = $shp11.Duplicate().Item(1)
$new15.Name = "Rectangle 14"
$new15.Left = 125.307400    # 1591404 EMU
$new15.Top  = 404.480790    # 5136906 EMU
$new15.TextFrame.TextRange.Text = "After Save"

# New "Before Change" box (accent5 style, like shp7) in the right column.
$new16 = $shp7.Duplicate().Item(1)
$new16.Name = "Rectangle 15"
$new16.Left = 560.220002    # 7114794 EMU
$new16.Top  = 170.740395    # 2168403 EMU
$new16.TextFrame.TextRange.Text = "Before Change"

# New "After Change" box (accent5 style, like shp9) in the right column.
$new18 = $shp9.Duplicate().Item(1)
$new18.Name = "Rectangle 17"
$new18.Left = 560.220002    # 7114794 EMU
$new18.Top  = 221.149056    # 2808593 EMU
$new18.TextFrame.TextRange.Text = "After Change"

# New "Before Collection Entity Initialize" box (accent2 style, like shp11) in the right column.
$new19 = $shp11.Duplicate().Item(1)
$new19.Name = "Rectangle 18"
$new19.Left = 560.219941    # 7114793 EMU
$new19.Top  = 361.749131    # 4594214 EMU
$new19.TextFrame.TextRange.Text = "Before Collection Entity Initialize"
$new19.TextFrame.TextRange.Font.Size = 14

# ---------------------------------------------------------------------------
# C) Re-label the old "After Save" box -> "Before Collection Entity Initialize"
#    (done last so the clones above still pick up the original "After Save"
#    run formatting).
# ---------------------------------------------------------------------------
$shp11.TextFrame.TextRange.Text = "Before Collection Entity Initialize"
$shp11.TextFrame.TextRange.Font.Size = 14
